$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore revision: update rule R30's "Integer max" value (cell C10) from 18 to 1
$ws.Range("C10").Value = 1
